# Atualizacao de bases das ligas, do dia: 11-04-2024 as 00:31
# Updates odds/match data rows for "Netherlands Eredivisie" sheet:
# rows 116/117, 168/169, 170/171 and 271/272 swap their match records
# (team names, scores, odds, etc.), plus independent odds refreshes on
# rows 266-274 (columns N/O/P/Q/R/S/U/V).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B116").Value = 6838439
$ws.Range("F116").Value = "PEC Zwolle"
$ws.Range("G116").Value = "RKC"
$ws.Range("H116").Value = 1
$ws.Range("I116").Value = 2
$ws.Range("J116").Value = "A"
$ws.Range("K116").Value = 1.909
$ws.Range("L116").Value = 3.6
$ws.Range("M116").Value = 3.75
$ws.Range("N116").Value = 2.05
$ws.Range("O116").Value = 3.6
$ws.Range("P116").Value = 3.4
$ws.Range("Q116").Value = -0.5
$ws.Range("R116").Value = 2.05
$ws.Range("S116").Value = 1.8
$ws.Range("T116").Value = 2.75
$ws.Range("U116").Value = 1.9
$ws.Range("V116").Value = 1.95
$ws.Range("W116").Value = -1
$ws.Range("Y116").Value = 2.4
$ws.Range("Z116").Value = -1
$ws.Range("AA116").Value = 0.8
$ws.Range("AB116").Value = 0.45
$ws.Range("AC116").Value = -0.5
$ws.Range("B117").Value = 6838440
$ws.Range("F117").Value = "Ajax"
$ws.Range("G117").Value = "Vitesse"
$ws.Range("H117").Value = 5
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = "H"
$ws.Range("K117").Value = 1.363
$ws.Range("L117").Value = 5
$ws.Range("M117").Value = 7.5
$ws.Range("N117").Value = 1.333
$ws.Range("O117").Value = 6
$ws.Range("P117").Value = 7.5
$ws.Range("Q117").Value = -1.5
$ws.Range("R117").Value = 1.825
$ws.Range("S117").Value = 2.025
$ws.Range("T117").Value = 3.5
$ws.Range("U117").Value = 1.95
$ws.Range("V117").Value = 1.9
$ws.Range("W117").Value = 0.333
$ws.Range("Y117").Value = -1
$ws.Range("Z117").Value = 0.825
$ws.Range("AA117").Value = -1
$ws.Range("AB117").Value = 0.95
$ws.Range("AC117").Value = -1
$ws.Range("B168").Value = 6838488
$ws.Range("F168").Value = "RKC"
$ws.Range("G168").Value = "Sparta Rotterdam"
$ws.Range("H168").Value = 1
$ws.Range("I168").Value = 1
$ws.Range("J168").Value = "D"
$ws.Range("K168").Value = 3.1
$ws.Range("L168").Value = 3.75
$ws.Range("M168").Value = 2.1
$ws.Range("N168").Value = 2.3
$ws.Range("O168").Value = 3.75
$ws.Range("P168").Value = 2.8
$ws.Range("Q168").Value = -0.25
$ws.Range("R168").Value = 2.025
$ws.Range("S168").Value = 1.825
$ws.Range("T168").Value = 2.5
$ws.Range("U168").Value = 1.975
$ws.Range("V168").Value = 1.875
$ws.Range("W168").Value = -1
$ws.Range("X168").Value = 2.75
$ws.Range("Z168").Value = -0.5
$ws.Range("AA168").Value = 0.4125
$ws.Range("AC168").Value = 0.875
$ws.Range("B169").Value = 6838489
$ws.Range("F169").Value = "PSV"
$ws.Range("G169").Value = "Almere City FC"
$ws.Range("H169").Value = 2
$ws.Range("I169").Value = 0
$ws.Range("J169").Value = "H"
$ws.Range("K169").Value = 1.1
$ws.Range("L169").Value = 9.5
$ws.Range("M169").Value = 19
$ws.Range("N169").Value = 1.1
$ws.Range("O169").Value = 10
$ws.Range("P169").Value = 23
$ws.Range("Q169").Value = -2.5
$ws.Range("R169").Value = 1.9
$ws.Range("S169").Value = 1.95
$ws.Range("T169").Value = 3.75
$ws.Range("U169").Value = 2.025
$ws.Range("V169").Value = 1.825
$ws.Range("W169").Value = 0.1000000000000001
$ws.Range("X169").Value = -1
$ws.Range("Z169").Value = -1
$ws.Range("AA169").Value = 0.95
$ws.Range("AC169").Value = 0.825
$ws.Range("B170").Value = 6838491
$ws.Range("F170").Value = "Heracles"
$ws.Range("G170").Value = "Ajax"
$ws.Range("H170").Value = 2
$ws.Range("I170").Value = 4
$ws.Range("J170").Value = "A"
$ws.Range("K170").Value = 7.5
$ws.Range("L170").Value = 5
$ws.Range("M170").Value = 1.363
$ws.Range("N170").Value = 4.75
$ws.Range("O170").Value = 4.5
$ws.Range("P170").Value = 1.6
$ws.Range("Q170").Value = 1
$ws.Range("R170").Value = 1.85
$ws.Range("S170").Value = 2
$ws.Range("T170").Value = 3.25
$ws.Range("U170").Value = 1.875
$ws.Range("V170").Value = 1.975
$ws.Range("X170").Value = -1
$ws.Range("Y170").Value = 0.6000000000000001
$ws.Range("Z170").Value = -1
$ws.Range("AA170").Value = 1
$ws.Range("AB170").Value = 0.875
$ws.Range("AC170").Value = -1
$ws.Range("B171").Value = 6838490
$ws.Range("F171").Value = "Excelsior"
$ws.Range("G171").Value = "FC Utrecht"
$ws.Range("H171").Value = 1
$ws.Range("I171").Value = 1
$ws.Range("J171").Value = "D"
$ws.Range("K171").Value = 2.7
$ws.Range("L171").Value = 3.6
$ws.Range("M171").Value = 2.375
$ws.Range("N171").Value = 3.2
$ws.Range("O171").Value = 3.6
$ws.Range("P171").Value = 2.15
$ws.Range("Q171").Value = 0.25
$ws.Range("R171").Value = 1.975
$ws.Range("S171").Value = 1.875
$ws.Range("T171").Value = 2.75
$ws.Range("U171").Value = 1.975
$ws.Range("V171").Value = 1.875
$ws.Range("X171").Value = 2.6
$ws.Range("Y171").Value = -1
$ws.Range("Z171").Value = 0.4875
$ws.Range("AA171").Value = -0.5
$ws.Range("AB171").Value = -1
$ws.Range("AC171").Value = 0.875
$ws.Range("N266").Value = 1.666
$ws.Range("O266").Value = 4.333
$ws.Range("P266").Value = 4.5
$ws.Range("R266").Value = 1.92
$ws.Range("S266").Value = 1.98
$ws.Range("N267").Value = 1.055
$ws.Range("O267").Value = 15
$ws.Range("P267").Value = 29
$ws.Range("R267").Value = 1.98
$ws.Range("S267").Value = 1.92
$ws.Range("U267").Value = 1.925
$ws.Range("V267").Value = 1.925
$ws.Range("N268").Value = 1.333
$ws.Range("O268").Value = 5.25
$ws.Range("R268").Value = 1.9
$ws.Range("S268").Value = 2
$ws.Range("U268").Value = 1.825
$ws.Range("V268").Value = 2.025
$ws.Range("N269").Value = 2.75
$ws.Range("P269").Value = 2.55
$ws.Range("R269").Value = 2.02
$ws.Range("S269").Value = 1.88
$ws.Range("N270").Value = 2.8
$ws.Range("O270").Value = 3.4
$ws.Range("P270").Value = 2.5
$ws.Range("Q270").Value = 0
$ws.Range("R270").Value = 2.08
$ws.Range("S270").Value = 1.82
$ws.Range("B271").Value = 6994878
$ws.Range("F271").Value = "Fortuna Sittard"
$ws.Range("G271").Value = "Feyenoord"
$ws.Range("K271").Value = 9
$ws.Range("L271").Value = 5.5
$ws.Range("M271").Value = 1.25
$ws.Range("N271").Value = 8
$ws.Range("O271").Value = 5.5
$ws.Range("P271").Value = 1.333
$ws.Range("Q271").Value = 1.5
$ws.Range("R271").Value = 1.9
$ws.Range("S271").Value = 2
$ws.Range("T271").Value = 3
$ws.Range("U271").Value = 2
$ws.Range("V271").Value = 1.85
$ws.Range("B272").Value = 6838573
$ws.Range("F272").Value = "FC Utrecht"
$ws.Range("G272").Value = "Go Ahead Eagles"
$ws.Range("K272").Value = 1.909
$ws.Range("L272").Value = 3.5
$ws.Range("M272").Value = 3.5
$ws.Range("N272").Value = 1.727
$ws.Range("O272").Value = 3.8
$ws.Range("P272").Value = 4.75
$ws.Range("Q272").Value = -0.75
$ws.Range("R272").Value = 1.99
$ws.Range("S272").Value = 1.91
$ws.Range("T272").Value = 2.75
$ws.Range("U272").Value = 1.85
$ws.Range("V272").Value = 2
$ws.Range("N273").Value = 2.7
$ws.Range("O273").Value = 3.4
$ws.Range("P273").Value = 2.6
$ws.Range("O274").Value = 5
$ws.Range("P274").Value = 5
$ws.Range("R274").Value = 1.92
$ws.Range("S274").Value = 1.98
